$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = '65.353.93'
$ws.Range("E2").Value = '  -6.34%  '
# Row 3
$ws.Range("D3").Value = '3.299.60'
$ws.Range("E3").Value = '  -7.14%  '
# Row 4
Set-TextValue "D4" '0.999'
$ws.Range("E4").Value = '  -0.18%  '
# Row 5
Set-TextValue "D5" '552.92'
$ws.Range("E5").Value = '  -5.94%  '
# Row 6
Set-TextValue "D6" '179.33'
$ws.Range("E6").Value = '  -8.96%  '
# Row 7
$ws.Range("E7").Value = '  +0.10%  '
# Row 8
$ws.Range("E8").Value = '  -4.20%  '
# Row 9
$ws.Range("D9").Value = '3.297.84'
$ws.Range("E9").Value = '  -6.84%  '
# Row 10
Set-TextValue "D10" '0.185'
$ws.Range("E10").Value = '  -12.23%  '
# Row 11
Set-TextValue "D11" '0.582'
$ws.Range("E11").Value = '  -7.14%  '
# Row 12
Set-TextValue "D12" '46.99'
$ws.Range("E12").Value = '  -10.98%  '
# Row 13
Set-TextValue "D13" '0.0000262'
$ws.Range("E13").Value = '  -9.20%  '
# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue "D14" '8.50'
$ws.Range("E14").Value = '  -8.05%  '
# Row 15
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '3.819.83'
$ws.Range("E15").Value = '  -7.39%  '
# Row 16
Set-TextValue "D16" '603.84'
$ws.Range("E16").Value = '  -8.86%  '
# Row 17
Set-TextValue "D17" '18.01'
$ws.Range("E17").Value = '  -2.46%  '
# Row 18
$ws.Range("D18").Value = '65.259.35'
$ws.Range("E18").Value = '  -6.32%  '
# Row 19
$ws.Range("E19").Value = '  -4.14%  '
# Row 20
$ws.Range("D20").Value = '3.289.67'
$ws.Range("E20").Value = '  -7.59%  '
# Row 21
Set-TextValue "D21" '11.35'
$ws.Range("E21").Value = '  -9.69%  '
# Row 22
$ws.Range("E22").Value = '  -6.83%  '
# Row 23
Set-TextValue "D23" '17.44'
$ws.Range("E23").Value = '  -3.64%  '
# Row 24
Set-TextValue "D24" '102.42'
$ws.Range("E24").Value = '  -2.41%  '
# Row 25
Set-TextValue "D25" '4.95'
$ws.Range("E25").Value = '  -7.66%  '
# Row 26
$ws.Range("E26").Value = '  -10.05%  '
# Row 27
Set-TextValue "D27" '5.97'
$ws.Range("E27").Value = '  -1.54%  '
# Row 28
$ws.Range("E28").Value = '  -9.16%  '
# Row 29
Set-TextValue "D29" '9.31'
$ws.Range("E29").Value = '  -8.84%  '
# Row 30
Set-TextValue "D30" '8.62'
$ws.Range("E30").Value = '  -10.25%  '
# Row 31
Set-TextValue "D31" '30.32'
$ws.Range("E31").Value = '  -9.18%  '
# Row 32
Set-TextValue "D32" '3.90'
$ws.Range("E32").Value = '  -11.38%  '
# Row 33
Set-TextValue "D33" '6.21'
$ws.Range("E33").Value = '  -8.99%  '
# Row 34
Set-TextValue "D34" '10.95'
$ws.Range("E34").Value = '  -7.08%  '
# Row 35
$ws.Range("D35").Value = '3.799.55'
$ws.Range("E35").Value = '  +1.21%  '
# Row 36
$ws.Range("E36").Value = '  -6.60%  '
# Row 37
Set-TextValue "D37" '530.06'
$ws.Range("E37").Value = '  +3.50%  '
# Row 38
Set-TextValue "D38" '1.00'
$ws.Range("E38").Value = '  +0.08%  '
# Row 39
Set-TextValue "D39" '56.00'
$ws.Range("E39").Value = '  -9.41%  '
# Row 40
$ws.Range("E40").Value = '  -9.60%  '
# Row 41
$ws.Range("D41").Value = '0.0₃0708'
$ws.Range("E41").Value = '  -13.57%  '
# Row 42
Set-TextValue "D42" '2.64'
$ws.Range("E42").Value = '  -10.45%  '
# Row 43
$ws.Range("E43").Value = '  -7.21%  '
# Row 44
$ws.Range("E44").Value = '  -9.56%  '
# Row 45
$ws.Range("E45").Value = '  -9.27%  '
# Row 46
Set-TextValue "D46" '3.18'
$ws.Range("E46").Value = '  +17.75%  '
# Row 47
Set-TextValue "D47" '3.23'
$ws.Range("E47").Value = '  -5.21%  '
# Row 48
Set-TextValue "D48" '0.0408'
$ws.Range("E48").Value = '  -10.52%  '
# Row 49
$ws.Range("E49").Value = '  -6.06%  '
# Row 50
Set-TextValue "D50" '2.57'
# Row 51
Set-TextValue "D51" '0.997'
$ws.Range("E51").Value = '  -0.36%  '
